$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D3").Value = "2016-01-20 03:28:24"
$wsZhCn.Range("G3").Value = "2016-01-20 03:29:10"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D3").Value = "2016-01-20 03:28:35"
$wsDeDe.Range("G3").Value = "2016-01-20 03:29:28"
